{"js": "// Word Office.js (Word JavaScript API) edit script.\n// Body of: async (context) => { ... }\n//\n// Net change (per the supplied diff):\n//  1. Right after the title paragraph (\"Play Dragon Born Slot for Free -\n//     Review 2021\", Heading1), insert a new Normal paragraph containing a\n//     bold \"Meta description\" run followed by a plain run with the rest of\n//     the meta description text (preceded by an empty run, matching the\n//     document's existing paragraph-run convention).\n//  2. Near the end of the document, delete the duplicate bold\n//     \"Play Dragon Born Slot for Free - Review 2021\" paragraph, and replace\n//     the text of the following italic paragraph with a new DALLE image\n//     prompt, keeping the italic formatting and leading empty run intact.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\n// ---- 1. Insert the \"Meta description\" paragraph after the title ----\n\nconst titlePara = paragraphs.items[0]; // \"Play Dragon Born Slot for Free - Review 2021\" (Heading1)\n\n// Create a new, empty paragraph right after the title.\nconst metaPara = titlePara.insertParagraph(\"\", \"After\");\nawait context.sync();\n\n// Use insertOoxml to give the new paragraph the exact run structure used\n// throughout this document (leading empty run, then the formatted runs),\n// and to make sure it picks up no paragraph style (i.e. \"Normal\").\nconst metaOoxml = `<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>\n<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">\n  <pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">\n    <pkg:xmlData>\n      <w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">\n        <w:body>\n          <w:p>\n            <w:r/>\n            <w:r>\n              <w:rPr><w:b/></w:rPr>\n              <w:t>Meta description</w:t>\n            </w:r>\n            <w:r>\n              <w:t>: Dragon Born is a Megaways slot game with up to 50 free spins. Read our review and play for free.</w:t>\n            </w:r>\n          </w:p>\n        </w:body>\n      </w:document>\n    </pkg:xmlData>\n  </pkg:part>\n</pkg:package>`;\n\nmetaPara.getRange(\"Whole\").insertOoxml(metaOoxml, \"Replace\");\nawait context.sync();\n\n// ---- 2. Replace the trailing duplicate-title / description paragraphs ----\n\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst items = paragraphs.items;\nconst count = items.length;\n// The last paragraph is the italic meta-description-like paragraph; the one\n// right before it is the duplicate bold title paragraph that must go away.\nconst duplicateTitlePara = items[count - 2];\nconst descriptionPara = items[count - 1];\n\nduplicateTitlePara.delete();\nawait context.sync();\n\nconst newDescription =\n  \"Create a Cartoon-Style Image Featuring a Happy Maya Warrior with Glasses for Dragon Born DALLE, we need your artistic skills to create an eye-catching feature image for Dragon Born. The image must fit the game's theme and feature a happy Maya warrior with glasses. We want to see a cartoon-style illustration that is vibrant, colorful, and energetic. The Maya warrior should be the focal point of the image, striking a dynamic pose and looking directly at the viewer with a big smile on their face. They should be holding a shield and a sword, ready for battle. The background should be set against a backdrop of a medieval castle, with flags flying in the wind. There should be a dragon flying in the sky, adding an element of danger and excitement to the image. The colors should be vibrant and bold, with plenty of contrast to make the image pop. The Maya warrior's outfit should be adorned with bright colors and intricate details, adding to the fantasy theme. Overall, we want an image that captures the spirit of Dragon Born: exciting, adventurous, and full of energy. We look forward to seeing your creative masterpiece!\";\n\ndescriptionPara.getRange(\"Whole\").insertText(newDescription, \"Replace\");\nawait context.sync();\n", "ps1": "# Word COM interop (PowerShell-style) edit script.\n# $word.ActiveDocument is already open as $d.\n#\n# Net change (per the supplied diff):\n#  1. Right after the title paragraph (\"Play Dragon Born Slot for Free -\n#     Review 2021\", Heading1), insert a new Normal paragraph containing a\n#     bold \"Meta description\" run followed by a plain run with the rest of\n#     the meta description text (preceded by an empty run, matching the\n#     document's existing paragraph-run convention).\n#  2. Near the end of the document, delete the duplicate bold\n#     \"Play Dragon Born Slot for Free - Review 2021\" paragraph, and replace\n#     the text of the following italic paragraph with a new DALLE image\n#     prompt, keeping the italic formatting and leading empty run intact.\n\n$d = $word.ActiveDocument\n\n# ---- 1. Insert the \"Meta description\" paragraph after the title ----\n\n$titleRange = $d.Paragraphs(1).Range\n$titleRange.InsertParagraphAfter()\n\n$metaParaRange = $d.Paragraphs(2).Range\n\n$metaOoxml = '<?xml version=\"1.0\" encoding=\"UTF-8\" standalone=\"yes\"?>' + \"`n\" `\n  + '<pkg:package xmlns:pkg=\"http://schemas.microsoft.com/office/2006/xmlPackage\">' + \"`n\" `\n  + '<pkg:part pkg:name=\"/word/document.xml\" pkg:contentType=\"application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml\">' + \"`n\" `\n  + '<pkg:xmlData>' + \"`n\" `\n  + '<w:document xmlns:w=\"http://schemas.openxmlformats.org/wordprocessingml/2006/main\">' `\n  + '<w:body>' `\n  + '<w:p>' `\n  + '<w:r/>' `\n  + '<w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r>' `\n  + '<w:r><w:t>: Dragon Born is a Megaways slot game with up to 50 free spins. Read our review and play for free.</w:t></w:r>' `\n  + '</w:p>' `\n  + '</w:body>' `\n  + '</w:document>' + \"`n\" `\n  + '</pkg:xmlData></pkg:part></pkg:package>'\n\n$metaParaRange.InsertXML($metaOoxml)\n\n# ---- 2. Replace the trailing duplicate-title / description paragraphs ----\n\n$count = $d.Paragraphs.Count\n$duplicateTitleRange = $d.Paragraphs($count - 1).Range\n$duplicateTitleRange.Delete()\n\n$newDescription = \"Create a Cartoon-Style Image Featuring a Happy Maya Warrior with Glasses for Dragon Born DALLE, we need your artistic skills to create an eye-catching feature image for Dragon Born. The image must fit the game's theme and feature a happy Maya warrior with glasses. We want to see a cartoon-style illustration that is vibrant, colorful, and energetic. The Maya warrior should be the focal point of the image, striking a dynamic pose and looking directly at the viewer with a big smile on their face. They should be holding a shield and a sword, ready for battle. The background should be set against a backdrop of a medieval castle, with flags flying in the wind. There should be a dragon flying in the sky, adding an element of danger and excitement to the image. The colors should be vibrant and bold, with plenty of contrast to make the image pop. The Maya warrior's outfit should be adorned with bright colors and intricate details, adding to the fantasy theme. Overall, we want an image that captures the spirit of Dragon Born: exciting, adventurous, and full of energy. We look forward to seeing your creative masterpiece!\"\n\n$descriptionRange = $d.Paragraphs($d.Paragraphs.Count).Range\n$descriptionRange.MoveEnd(1, -1)   # wdCharacter = 1; exclude the trailing paragraph mark\n$descriptionRange.Text = $newDescription\n"}
